# update format for output
# Rewrites the "shock" / "extreme_level" columns from raw numeric values to
# formatted text labels, renumbers the tail of the table (old M15-M19 rows
# are dropped, old M20-M25 rows shift up and get new labels/values), and
# shrinks the used range from A1:D31 down to A1:D21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $text) {
    # A leading apostrophe forces Excel to store the value as literal text
    # instead of auto-coercing percentage-looking strings (e.g. "0.1 %")
    # into numbers; ClearFormats then drops the resulting quote-prefix
    # style so the cell ends up with no explicit style, matching the rest
    # of the data cells on this sheet.
    $cell = $sheet.Range($addr)
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

function Clear-CellText($sheet, $addr) {
    # Assigning "" outright makes Excel drop the cell back to an untyped
    # blank (Number/null) instead of an empty text cell. A lone leading
    # apostrophe forces Excel's literal-text path, landing on an empty
    # string (matching the other blank "extreme_level" cells on this
    # sheet, which are empty inline strings); ClearFormats then removes
    # the quote-prefix style so no stray "s" attribute is left behind.
    $cell = $sheet.Range($addr)
    $cell.Value = "'"
    $cell.ClearFormats()
}

# --- Rows 2-10 (M1..M9): shock column becomes "0.1 %", extreme_level stays blank ---
Set-TextValue $ws "C2"  "0.1 %"
Set-TextValue $ws "C3"  "0.1 %"
Set-TextValue $ws "C4"  "0.1 %"
Set-TextValue $ws "C5"  "0.1 %"
Set-TextValue $ws "C6"  "0.1 %"
Set-TextValue $ws "C7"  "0.1 %"
Set-TextValue $ws "C8"  "0.1 %"
Set-TextValue $ws "C9"  "0.1 %"
Set-TextValue $ws "C10" "0.1 %"

# --- Row 11 (M10): shock becomes "1.4 %" ---
Set-TextValue $ws "C11" "1.4 %"

# --- Rows 12-15 (M11..M14): shock + extreme_level become text labels ---
Set-TextValue $ws "C12" "13% max"
Set-TextValue $ws "D12" "(+972 ppts)"

Set-TextValue $ws "C13" "13 peak"
Set-TextValue $ws "D13" "(+990 ppts)"

Set-TextValue $ws "C14" "13% peak"
Set-TextValue $ws "D14" "(+100800 bps)"

Set-TextValue $ws "C15" "13% peak"
Set-TextValue $ws "D15" "(+102600 bps)"

# --- Row 16: was M15/MMM15 -> becomes M20/MMM20 with new values ---
Set-TextValue $ws "A16" "M20"
Set-TextValue $ws "B16" "MMM20"
Set-TextValue $ws "C16" "0.4 %"
Set-TextValue $ws "D16" "0.4 %"

# --- Row 17: was M16/MMM16 -> becomes M21/MMM21 with new values ---
Set-TextValue $ws "A17" "M21"
Set-TextValue $ws "B17" "MMM21"
Set-TextValue $ws "C17" "13% peak"
Set-TextValue $ws "D17" "(+115200 bps)"

# --- Row 18: was M17/MMM17 -> becomes M22/MMM22 with new values ---
Set-TextValue $ws "A18" "M22"
Set-TextValue $ws "B18" "MMM22"
Set-TextValue $ws "C18" "13% peak"
Set-TextValue $ws "D18" "(+117000 bps)"

# --- Row 19: was M18/MMM18 -> becomes M23/MMM23 with new values ---
Set-TextValue $ws "A19" "M23"
Set-TextValue $ws "B19" "MMM23"
Set-TextValue $ws "C19" "13% peak"
Set-TextValue $ws "D19" "(+118800 bps)"

# --- Row 20: was M19/MMM19 -> becomes M24/MMM24 with new values ---
Set-TextValue $ws "A20" "M24"
Set-TextValue $ws "B20" "MMM24"
Set-TextValue $ws "C20" "0.1 %"
Clear-CellText  $ws "D20"

# --- Row 21: was M20/MMM20 -> becomes M25/MMM25 with new values ---
Set-TextValue $ws "A21" "M25"
Set-TextValue $ws "B21" "MMM25"
Set-TextValue $ws "C21" "13% peak"
Set-TextValue $ws "D21" "(+122400 bps)"

# --- Drop the old rows 22-31 (M21..M30) entirely; sheet now ends at row 21 ---
$ws.Range("A22:D31").ClearContents()
